$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 86: Scheduler (column T) cleared to empty string (text for next run not yet recorded)
$ws.Cells.Item(86, 20).Value = ""

# New training run rows 87-96 appended to the log
$newRows = @(
    @{ "A" = "2024-1-5 11:14:18"; "B" = 20; "C" = 64; "D" = 0.01; "E" = "ADAM"; "F" = "CEL"; "G" = 13.2; "H" = 32; "I" = 1.5145; "J" = 1.1936; "K" = 51.8792; "L" = "FER2013"; "M" = "cuda:0"; "N" = 4; "O" = 2; "P" = "Stationær"; "Q" = 263.1; "R" = 0; "S" = 0.005; "T" = ""; "U" = 0 },
    @{ "A" = "2024-1-5 11:16:24"; "B" = 20; "C" = 64; "D" = 0.01; "E" = "ADAM"; "F" = "CEL"; "G" = 14.2; "H" = 32; "I" = 1.272; "J" = 1.1194; "K" = 54.5926; "L" = "FER2013"; "M" = "cuda:0"; "N" = 4; "O" = 2; "P" = "Stationær"; "Q" = 284.1; "R" = 0; "S" = 0.005; "T" = ""; "U" = 0 },
    @{ "A" = "2024-1-5 11:16:38"; "B" = 20; "C" = 64; "D" = 0.01; "E" = "ADAM"; "F" = "CEL"; "G" = 14.3; "H" = 32; "I" = 1.2186; "J" = 1.2186; "K" = 53.1332; "L" = "FER2013"; "M" = "cuda:0"; "N" = 4; "O" = 2; "P" = "Stationær"; "Q" = 285.3; "R" = 0; "S" = 0.005; "T" = ""; "U" = 0 },
    @{ "A" = "2024-1-5 11:16:53"; "B" = 20; "C" = 64; "D" = 0.01; "E" = "ADAM"; "F" = "CEL"; "G" = 14.3; "H" = 32; "I" = 1.4361; "J" = 1.1578; "K" = 51.2871; "L" = "FER2013"; "M" = "cuda:0"; "N" = 4; "O" = 2; "P" = "Stationær"; "Q" = 286.5; "R" = 0; "S" = 0.005; "T" = ""; "U" = 0 },
    @{ "A" = "2024-1-5 11:17:6"; "B" = 20; "C" = 64; "D" = 0.01; "E" = "ADAM"; "F" = "CEL"; "G" = 14.4; "H" = 32; "I" = 0.9941; "J" = 0.9941; "K" = 53.579; "L" = "FER2013"; "M" = "cuda:0"; "N" = 4; "O" = 2; "P" = "Stationær"; "Q" = 287.5; "R" = 0; "S" = 0.005; "T" = ""; "U" = 0 },
    @{ "A" = "2024-1-5 11:17:14"; "B" = 20; "C" = 64; "D" = 0.01; "E" = "ADAM"; "F" = "CEL"; "G" = 14.3; "H" = 32; "I" = 1.4776; "J" = 1.1102; "K" = 53.6905; "L" = "FER2013"; "M" = "cuda:0"; "N" = 4; "O" = 2; "P" = "Stationær"; "Q" = 286.8; "R" = 0; "S" = 0.005; "T" = ""; "U" = 0 },
    @{ "A" = "2024-1-5 11:17:25"; "B" = 20; "C" = 64; "D" = 0.01; "E" = "ADAM"; "F" = "CEL"; "G" = 14.4; "H" = 32; "I" = 1.1259; "J" = 1.1253; "K" = 53.1436; "L" = "FER2013"; "M" = "cuda:0"; "N" = 4; "O" = 2; "P" = "Stationær"; "Q" = 288.8; "R" = 0; "S" = 0.005; "T" = ""; "U" = 0 },
    @{ "A" = "2024-1-5 11:17:34"; "B" = 20; "C" = 64; "D" = 0.01; "E" = "ADAM"; "F" = "CEL"; "G" = 14.3; "H" = 32; "I" = 1.1571; "J" = 1.1571; "K" = 52.8545; "L" = "FER2013"; "M" = "cuda:0"; "N" = 4; "O" = 2; "P" = "Stationær"; "Q" = 286.1; "R" = 0; "S" = 0.005; "T" = ""; "U" = 0 },
    @{ "A" = "2024-1-5 11:17:39"; "B" = 20; "C" = 64; "D" = 0.01; "E" = "ADAM"; "F" = "CEL"; "G" = 14.2; "H" = 32; "I" = 1.384; "J" = 1.2193; "K" = 51.3915; "L" = "FER2013"; "M" = "cuda:0"; "N" = 4; "O" = 2; "P" = "Stationær"; "Q" = 284.7; "R" = 0; "S" = 0.005; "T" = ""; "U" = 0 },
    @{ "A" = "2024-1-5 11:17:51"; "B" = 20; "C" = 64; "D" = 0.01; "E" = "ADAM"; "F" = "CEL"; "G" = 14.1; "H" = 32; "I" = 1.1389; "J" = 1.1389; "K" = 52.6908; "L" = "FER2013"; "M" = "cuda:0"; "N" = 4; "O" = 2; "P" = "Stationær"; "Q" = 282; "R" = 0; "S" = 0.005; "T" = "None"; "U" = 0 }
)

$rowNum = 87
foreach ($r in $newRows) {
    $ws.Cells.Item($rowNum, 1).Value = $r["A"]
    $ws.Cells.Item($rowNum, 2).Value = $r["B"]
    $ws.Cells.Item($rowNum, 3).Value = $r["C"]
    $ws.Cells.Item($rowNum, 4).Value = $r["D"]
    $ws.Cells.Item($rowNum, 5).Value = $r["E"]
    $ws.Cells.Item($rowNum, 6).Value = $r["F"]
    $ws.Cells.Item($rowNum, 7).Value = $r["G"]
    $ws.Cells.Item($rowNum, 8).Value = $r["H"]
    $ws.Cells.Item($rowNum, 9).Value = $r["I"]
    $ws.Cells.Item($rowNum, 10).Value = $r["J"]
    $ws.Cells.Item($rowNum, 11).Value = $r["K"]
    $ws.Cells.Item($rowNum, 12).Value = $r["L"]
    $ws.Cells.Item($rowNum, 13).Value = $r["M"]
    $ws.Cells.Item($rowNum, 14).Value = $r["N"]
    $ws.Cells.Item($rowNum, 15).Value = $r["O"]
    $ws.Cells.Item($rowNum, 16).Value = $r["P"]
    $ws.Cells.Item($rowNum, 17).Value = $r["Q"]
    $ws.Cells.Item($rowNum, 18).Value = $r["R"]
    $ws.Cells.Item($rowNum, 19).Value = $r["S"]
    $ws.Cells.Item($rowNum, 20).Value = $r["T"]
    $ws.Cells.Item($rowNum, 21).Value = $r["U"]
    $rowNum++
}